$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("D3").Value = "Order"
$ws.Range("F3").Value = "TOP3 "
$ws.Range("E3").Value = "MUSEUM"
$ws.Range("G3").Value = "GO/NOT"

# D column: running order/count of each choice (entered row-by-row, like a
# user typing the formula into D4 then re-typing/filling each subsequent row)
for ($r = 4; $r -le 28; $r++) {
    $ws.Range("D$r").Formula = "=COUNTIF(`$C`$4:C$r,C$r)"
}

# E column: does this row match the currently selected museum?
$ws.Range("E4").Formula = "=`$C`$4:`$C`$28=`$I`$4"
$ws.Range("E5:E28").Formula = "=`$C`$4:`$C`$28=`$I`$4"

# F column: is this row within the capacity limit for the chosen museum?
$ws.Range("F4").Formula = "=D4:D4<=VLOOKUP(I`$10,I`$4:J`$7,2,FALSE)"
$ws.Range("F5:F28").Formula = "=D5:D5<=VLOOKUP(I`$10,I`$4:J`$7,2,FALSE)"

# G column: go / not go flag (match AND within capacity)
$ws.Range("G4").Formula = "=E4*F4"
$ws.Range("G5:G28").Formula = "=E5*F5"

# Column width adjustments to fit the new/edited columns
$ws.Columns.Item(1).ColumnWidth = 32.5703125
$ws.Columns.Item(2).ColumnWidth = 9
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Range("D1:G1").ColumnWidth = 9.28515625
$ws.Columns.Item(9).ColumnWidth = 13.140625
$ws.Columns.Item(12).ColumnWidth = 13.42578125

# Update the active selection to match the final cursor position
$ws.Range("J16").Select()
